$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "'299.10"
$ws.Range("E2").Value = "'2.01%"
$ws.Range("D3").Value = "'42.19"
$ws.Range("E3").Value = "'4.22%"
$ws.Range("D4").Value = "'5.016"
$ws.Range("E4").Value = "'0.11%"
$ws.Range("D5").Value = "'0.07546"
$ws.Range("E5").Value = "'2.54%"
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").Value = "'4.361"
$ws.Range("E6").Value = "'1.60%"
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").Value = "'1.604"
$ws.Range("E7").Value = "'2.32%"
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = "'0.9365"
$ws.Range("E8").Value = "'1.18%"
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").Value = "'2.389"
$ws.Range("E9").Value = "'1.58%"
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = "'0.1188"
$ws.Range("E10").Value = "'1.47%"
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = "'0.1839"
$ws.Range("E11").Value = "'2.74%"
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = "'0.09065"
$ws.Range("E12").Value = "'3.86%"
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = "'0.04161"
$ws.Range("E13").Value = "'-5.04%"
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = "'0.1047"
$ws.Range("E14").Value = "'-0.71%"
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = "'0.001294"
$ws.Range("E15").Value = "'2.19%"
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = "'0.005796"
$ws.Range("E16").Value = "'-2.02%"
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = "'3.340"
$ws.Range("E17").Value = "'-0.27%"
$ws.Range("D18").Value = "'0.3335"
$ws.Range("E18").Value = "'0.88%"
$ws.Range("D19").Value = "'8.334"
$ws.Range("E19").Value = "'6.43%"
$ws.Range("D20").Value = "'0.1399"
$ws.Range("E20").Value = "'0.72%"
$ws.Range("D21").Value = "'0.3100"
$ws.Range("E21").Value = "'11.92%"
$ws.Range("D22").Value = "'0.04083"
$ws.Range("E22").Value = "'4.25%"
$ws.Range("E23").Value = "'0.26%"
$ws.Range("D24").Value = "'0.003896"
$ws.Range("E24").Value = "'5.98%"
$ws.Range("D25").Value = "'0.0001299"
$ws.Range("E25").Value = "'8.26%"
$ws.Range("D38").Value = "'0.02408"
$ws.Range("E38").Value = "'2.74%"
$ws.Range("D39").Value = "'0.05233"
$ws.Range("E39").Value = "'2.62%"
$ws.Range("D40").Value = "'0.006752"
$ws.Range("E40").Value = "'21.83%"
$ws.Range("D41").Value = "'0.007745"
$ws.Range("E41").Value = "'-1.49%"
$ws.Range("D42").Value = "'0.1326"
$ws.Range("E42").Value = "'2.61%"
$ws.Range("D43").Value = "'0.007387"
$ws.Range("E43").Value = "'0.17%"
$ws.Range("D44").Value = "'0.007121"
$ws.Range("E44").Value = "'-11.62%"
$ws.Range("D45").Value = "'0.2997"
$ws.Range("E45").Value = "'2.84%"
$ws.Range("D46").Value = "'0.00006237"
$ws.Range("E46").Value = "'0.12%"
$ws.Range("E47").Value = "'0.00%"
$ws.Range("D48").Value = "'0.04574"
$ws.Range("E48").Value = "'-5.47%"
$ws.Range("E50").Value = "'0.00%"
$ws.Range("E51").Value = "'0.00%"
